# Musician Finder Specification Document - "Chat Area" slide (slide 12)
# 1) Shrink/move the "Rectangle 4" placeholder rectangle (top moves down, height shrinks).
# 2) Add a new "Name" label textbox above it (the chat-area header).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- 1. Resize/move the existing "Rectangle 4" shape ---
$rect4 = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Rectangle 4") {
        $rect4 = $s.Shapes.Item($i)
        break
    }
}

$rect4.Left   = 341.25
$rect4.Top    = 144.97937007874015
$rect4.Width  = 235.5
$rect4.Height = 272.7707086614173

# --- 2. Add the new "Name" textbox ---
# The slide already "used up" shape ids 3 and 9 earlier in its edit history (they
# belonged to shapes that were later removed), so the next two calls to
# AddTextbox on this slide reproduce that numbering before we create the one we
# actually want to keep, landing it on Id=10 / Name="TextBox 9" just like the
# authored deck.
$tmp1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tmp2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tmp1.Delete()

$nameBox = $s.Shapes.AddTextbox(1, 351.75, 104.85433070866142, 218.99992125984252, 29.081259842519685)

$tmp2.Delete()

$nameBox.TextFrame.WordWrap = -1
$nameBox.TextFrame.AutoSize = 1
$nameBox.Fill.Visible = 0
$nameBox.TextFrame.TextRange.Text = "Name"
